# Append two new paragraphs (each preceded by a blank paragraph) to the end
# of the memo body, right before the closing section properties, matching
# the "Child Chore page + cashout history complete" progress-report update.

$d = $word.ActiveDocument

$paragraph1 = 'For the next week I want to get a skeleton of the settings page done with little to no functionality, because it is a low priority for me. Other than that, I want to get the controller started and figure out how local storage works with Xamarin to store the accounts associated with the app. After that I want to start getting the child side GUI flushed out. If I get far enough into that this week, getting actual account swapping to work will be the next step.'

$paragraph2 = 'Only stopping issue that I might see is that I’ve somehow already used about 50% of my free azure credits testing and doing the original implementation, so if I run out I may have to halt progress to find another database hosting site to continue going.'

# Locate the last paragraph in the document (the one ending in "...parent account.")
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $lastPara.Range
$r.Collapse(0)

# Blank line, new paragraph of text, blank line, new paragraph of text --
# mirrors the structure added in the diff (empty <w:p/> spacer paragraphs
# between the two new content paragraphs).
$r.InsertAfter("`r`r" + $paragraph1 + "`r`r" + $paragraph2)
